$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value, applied identically to both
# the "展览" (Exhibition) and "全部类型" (All types) sheets.
$updates = @{
    4  = 2938
    5  = 206
    8  = 1622
    11 = 347
    24 = 126
    27 = 1949
    31 = 154
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
